$wb = $excel.ActiveWorkbook

# --- Add the two new sheets at the end, in order: optimun, solution ---
$sheets = $wb.Worksheets
$lastSheet = $sheets.Item($sheets.Count)
$optimun = $sheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$optimun.Name = "optimun"

$solution = $sheets.Add([System.Reflection.Missing]::Value, $optimun)
$solution.Name = "solution"
